# Primera edicion de datos
# Updates the "Direccion" (F) column with real street addresses in place of
# placeholder/junk text, and fixes the "EstalCivil" (J) column values that
# read "Arrejuntado" to the correct "Conviviente".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (Direccion) - rows 2..15
$direcciones = @{
    2  = "120 Pheasant Drive Venice, FL 34293"
    3  = "593 SW. Center Dr. Cumberland, RI 02864"
    4  = "56 Warren St. Londonderry, NH 03053"
    5  = "8398 Prairie Rd. South Richmond Hill, NY 11419"
    6  = "45 East Orange Ave. Oak Creek, WI 53154"
    7  = "310 Lower River Street Revere, MA 02151"
    8  = "793 Cypress Street Loveland, OH 45140"
    9  = "44 Peg Shop St. Pottstown, PA 19464"
    10 = "7276 Hillcrest Rd. Staten Island, NY 10301"
    11 = "86 Court Lane South Portland, ME 04106"
    12 = "343 Indian Spring StreetArlington, MA 02474"
    13 = "7109 SW. Ann Circle Chevy Chase, MD 20815"
    14 = "581 Garfield St. Aiken, SC 29803"
    15 = "8463 W. Westport Road Minneapolis, MN 55406"
}

foreach ($row in $direcciones.Keys) {
    $ws.Cells.Item($row, 6).Value = $direcciones[$row]
}

# Column J (EstalCivil) - rows 7, 10, 13: "Arrejuntado" -> "Conviviente"
$estadoRows = @(7, 10, 13)
foreach ($row in $estadoRows) {
    $ws.Cells.Item($row, 10).Value = "Conviviente"
}

# Clear the lingering K2 selection left over from the previous session.
$ws.Range("A1").Select()
